$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing data rows 2-60 with new GDP per Capita values
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "1156"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "1170"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "1189"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1240"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "1245"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "1246"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "1278"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "1253"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "1298"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "1307"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "1326"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "1309"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "1302"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "1298"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1291"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "1290"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "1272"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1277"
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "1290"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "1278"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "1272"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "1237"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "1007"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "1011"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "1039"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "1074"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "1105"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "1022"
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "1070"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "1023"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "1019"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "1144"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "1270"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "1347"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "1337"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "1304"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "1344"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "1211"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "1101"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "999"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "963"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "881.170438588735"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "843.875348899883"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "578.402744031048"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "428.424558289543"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "632.940396261985"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "600.175264462535"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "570.598118955881"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "545.038752668709"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "518.657866732878"
$c.ClearFormats()
$c = $ws.Range("E52")
$c.NumberFormat = "@"
$c.Value = "502.372745598633"
$c.ClearFormats()
$c = $ws.Range("E53")
$c.NumberFormat = "@"
$c.Value = "489.681996907494"
$c.ClearFormats()
$c = $ws.Range("E54")
$c.NumberFormat = "@"
$c.Value = "796.816563246003"
$c.ClearFormats()
$c = $ws.Range("E55")
$c.NumberFormat = "@"
$c.Value = "842.805165479659"
$c.ClearFormats()
$c = $ws.Range("E56")
$c.NumberFormat = "@"
$c.Value = "869.039272786813"
$c.ClearFormats()
$c = $ws.Range("E57")
$c.NumberFormat = "@"
$c.Value = "964.4081099473"
$c.ClearFormats()
$c = $ws.Range("E58")
$c.NumberFormat = "@"
$c.Value = "1057.09655092876"
$c.ClearFormats()
$c = $ws.Range("E59")
$c.NumberFormat = "@"
$c.Value = "1259.99674275737"
$c.ClearFormats()
$c = $ws.Range("E60")
$c.NumberFormat = "@"
$c.Value = "1319.60736728553"
$c.ClearFormats()

# Add new rows 61-68 for years 2009-2016
$ws.Cells.Item(61, 1).Value = 4
$ws.Cells.Item(61, 2).Value = "Afghanistan"
$ws.Cells.Item(61, 3).Value = "GDP per Capita"
$ws.Cells.Item(61, 4).Value = 2009
$c = $ws.Cells.Item(61, 5)
$c.NumberFormat = "@"
$c.Value = "1557.32063657228"
$c.ClearFormats()
$ws.Cells.Item(62, 1).Value = 4
$ws.Cells.Item(62, 2).Value = "Afghanistan"
$ws.Cells.Item(62, 3).Value = "GDP per Capita"
$ws.Cells.Item(62, 4).Value = 2010
$c = $ws.Cells.Item(62, 5)
$c.NumberFormat = "@"
$c.Value = "1627.67163410066"
$c.ClearFormats()
$ws.Cells.Item(63, 1).Value = 4
$ws.Cells.Item(63, 2).Value = "Afghanistan"
$ws.Cells.Item(63, 3).Value = "GDP per Capita"
$ws.Cells.Item(63, 4).Value = 2011
$c = $ws.Cells.Item(63, 5)
$c.NumberFormat = "@"
$c.Value = "1792"
$c.ClearFormats()
$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Afghanistan"
$ws.Cells.Item(64, 3).Value = "GDP per Capita"
$ws.Cells.Item(64, 4).Value = 2012
$c = $ws.Cells.Item(64, 5)
$c.NumberFormat = "@"
$c.Value = "1945"
$c.ClearFormats()
$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = "Afghanistan"
$ws.Cells.Item(65, 3).Value = "GDP per Capita"
$ws.Cells.Item(65, 4).Value = 2013
$c = $ws.Cells.Item(65, 5)
$c.NumberFormat = "@"
$c.Value = "2025"
$c.ClearFormats()
$ws.Cells.Item(66, 1).Value = 4
$ws.Cells.Item(66, 2).Value = "Afghanistan"
$ws.Cells.Item(66, 3).Value = "GDP per Capita"
$ws.Cells.Item(66, 4).Value = 2014
$c = $ws.Cells.Item(66, 5)
$c.NumberFormat = "@"
$c.Value = "2022"
$c.ClearFormats()
$ws.Cells.Item(67, 1).Value = 4
$ws.Cells.Item(67, 2).Value = "Afghanistan"
$ws.Cells.Item(67, 3).Value = "GDP per Capita"
$ws.Cells.Item(67, 4).Value = 2015
$c = $ws.Cells.Item(67, 5)
$c.NumberFormat = "@"
$c.Value = "1928"
$c.ClearFormats()
$ws.Cells.Item(68, 1).Value = 4
$ws.Cells.Item(68, 2).Value = "Afghanistan"
$ws.Cells.Item(68, 3).Value = "GDP per Capita"
$ws.Cells.Item(68, 4).Value = 2016
$c = $ws.Cells.Item(68, 5)
$c.NumberFormat = "@"
$c.Value = "1929"
$c.ClearFormats()
